$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 39.27143501396318
$ws.Range("G2").Value = 3.000336175386721
$ws.Range("H2").Value = 112.5454926273685
$ws.Range("I2").Value = 4.263868887090577
$ws.Range("J2").Value = 4.262443736231989
$ws.Range("K2").Value = 4.265761030084986
$ws.Range("L2").Value = 0.2393417146442213
$ws.Range("M2").Value = 0.2324561146610944
$ws.Range("N2").Value = 0.2531750388970516

$ws.Range("F3").Value = 0.2496588112822442
$ws.Range("G3").Value = 0.2494438598988644
$ws.Range("H3").Value = 0.2498914445169295
$ws.Range("I3").Value = 0.2235324113426677
$ws.Range("J3").Value = 0.2233345540825502
$ws.Range("K3").Value = 0.2237450427264582
$ws.Range("L3").Value = 0.2472445740366078
$ws.Range("M3").Value = 0.2470310257008661
$ws.Range("N3").Value = 0.247475808417009

$ws.Range("F4").Value = 39.52109382524542
$ws.Range("G4").Value = 3.249780035285586
$ws.Range("H4").Value = 112.7953840718854
$ws.Range("I4").Value = 4.487401298433246
$ws.Range("J4").Value = 4.485778290314539
$ws.Range("K4").Value = 4.489506072811445
$ws.Range("L4").Value = 0.486586288680829
$ws.Range("M4").Value = 0.4794871403619605
$ws.Range("N4").Value = 0.5006508473140605
